# Auto-generated edit script applying scheduled market-data refresh
# to the profit-analysis tables on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 111112520
$ws.Range("I100").Value = 125000960
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 125000960
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -125000419
$ws.Range("N100").Value = -6082
$ws.Range("H116").Value = 4950.375
$ws.Range("I116").Value = 2150
$ws.Range("K116").Value = 2150
$ws.Range("M116").Value = 1292
$ws.Range("H129").Value = 1125.1904
$ws.Range("I129").Value = 363.33334
$ws.Range("J129").Value = 1252.1666
$ws.Range("K129").Value = 1090.00002
$ws.Range("L129").Value = 3756.4998
$ws.Range("M129").Value = 3909.99998
$ws.Range("N129").Value = -13756.4998
$ws.Range("H132").Value = 2955.025
$ws.Range("I132").Value = 2935.6453
$ws.Range("J132").Value = 3021.7778
$ws.Range("K132").Value = 8806.9359
$ws.Range("L132").Value = 9065.3334
$ws.Range("M132").Value = -6276.9359
$ws.Range("N132").Value = -14125.3334
$ws.Range("H138").Value = 71434070
$ws.Range("J138").Value = 4746.6665
$ws.Range("L138").Value = 14239.9995
$ws.Range("N138").Value = -24519.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3116.6667
$ws.Range("I2").Value = 2740
$ws.Range("K2").Value = 2740
$ws.Range("M2").Value = -2627
$ws.Range("H32").Value = 5544.72
$ws.Range("I32").Value = 4483.1
$ws.Range("J32").Value = 9791.200000000001
$ws.Range("K32").Value = 4483.1
$ws.Range("L32").Value = 9791.200000000001
$ws.Range("M32").Value = -4196.1
$ws.Range("N32").Value = -10365.2
$ws.Range("H116").Value = 3116.6667
$ws.Range("I116").Value = 2740
$ws.Range("K116").Value = 2740
$ws.Range("M116").Value = -446
$ws.Range("H132").Value = 16876.146
$ws.Range("I132").Value = 1901.3334
$ws.Range("J132").Value = 74636.14
$ws.Range("K132").Value = 5704.0002
$ws.Range("L132").Value = 223908.42
$ws.Range("M132").Value = -3174.0002
$ws.Range("N132").Value = -228968.42

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3116.6667
$ws.Range("I3").Value = 2740
$ws.Range("K3").Value = 2740
$ws.Range("M3").Value = -2626
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H105").Value = 2944380
$ws.Range("I105").Value = 3441.818
$ws.Range("K105").Value = 3441.818
$ws.Range("M105").Value = -1694.818
$ws.Range("H107").Value = 735.73334
$ws.Range("I107").Value = 681.1429000000001
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 681.1429000000001
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1238.8571
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 550.1539
$ws.Range("I22").Value = 515
$ws.Range("J22").Value = 667.3333
$ws.Range("K22").Value = 515
$ws.Range("L22").Value = 667.3333
$ws.Range("M22").Value = -165
$ws.Range("N22").Value = -1367.3333
$ws.Range("H31").Value = 3388.6155
$ws.Range("I31").Value = 3210.9092
$ws.Range("K31").Value = 3210.9092
$ws.Range("M31").Value = -2915.9092
$ws.Range("H34").Value = 3388.6155
$ws.Range("I34").Value = 3210.9092
$ws.Range("K34").Value = 3210.9092
$ws.Range("M34").Value = -3008.9092
$ws.Range("H94").Value = 2089.6924
$ws.Range("I94").Value = 572.6667
$ws.Range("J94").Value = 3390
$ws.Range("K94").Value = 572.6667
$ws.Range("L94").Value = 3390
$ws.Range("M94").Value = -121.6667
$ws.Range("N94").Value = -4292
$ws.Range("H107").Value = 1357.7037
$ws.Range("I107").Value = 1117.1
$ws.Range("J107").Value = 1499.2354
$ws.Range("K107").Value = 1117.1
$ws.Range("L107").Value = 1499.2354
$ws.Range("M107").Value = 802.9000000000001
$ws.Range("N107").Value = -5339.2354
$ws.Range("H122").Value = 1975
$ws.Range("J122").Value = 1583.3334
$ws.Range("L122").Value = 4750.0002
$ws.Range("N122").Value = -9650.0002
$ws.Range("H134").Value = 1208.9546
$ws.Range("I134").Value = 1092.2
$ws.Range("J134").Value = 1459.1428
$ws.Range("K134").Value = 3276.6
$ws.Range("L134").Value = 4377.428400000001
$ws.Range("M134").Value = -741.6000000000004
$ws.Range("N134").Value = -9447.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 697.58
$ws.Range("J131").Value = 717.6129
$ws.Range("L131").Value = 2152.8387
$ws.Range("N131").Value = -12232.8387
$ws.Range("H132").Value = 568
$ws.Range("I132").Value = 568
$ws.Range("K132").Value = 5112
$ws.Range("M132").Value = -2582
$ws.Range("H140").Value = 1765.4445
$ws.Range("J140").Value = 2999.625
$ws.Range("L140").Value = 8998.875
$ws.Range("N140").Value = -19358.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H70").Value = 3133483.8
$ws.Range("I70").Value = 3871.25
$ws.Range("J70").Value = 5219892
$ws.Range("K70").Value = 3871.25
$ws.Range("L70").Value = 5219892
$ws.Range("M70").Value = -3601.25
$ws.Range("N70").Value = -5220432
$ws.Range("H73").Value = 3133483.8
$ws.Range("I73").Value = 3871.25
$ws.Range("J73").Value = 5219892
$ws.Range("K73").Value = 3871.25
$ws.Range("L73").Value = 5219892
$ws.Range("M73").Value = -2935.25
$ws.Range("N73").Value = -5221764
$ws.Range("H97").Value = 1767.5555
$ws.Range("I97").Value = 1947.9166
$ws.Range("J97").Value = 1406.8334
$ws.Range("K97").Value = 1947.9166
$ws.Range("L97").Value = 1406.8334
$ws.Range("M97").Value = -1451.9166
$ws.Range("N97").Value = -2398.8334
$ws.Range("H113").Value = 2309.05
$ws.Range("I113").Value = 1782.1538
$ws.Range("K113").Value = 1782.1538
$ws.Range("M113").Value = 387.8462
$ws.Range("H132").Value = 20657.297
$ws.Range("I132").Value = 1831.8235
$ws.Range("K132").Value = 5495.470499999999
$ws.Range("M132").Value = -2965.470499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 604388.9
$ws.Range("I132").Value = 804705.7
$ws.Range("K132").Value = 2414117.1
$ws.Range("M132").Value = -2411587.1
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -34860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19853380
$ws.Range("I136").Value = 26469148
$ws.Range("J136").Value = 6077.231
$ws.Range("K136").Value = 79407444
$ws.Range("L136").Value = 18231.693
$ws.Range("M136").Value = -79404894
$ws.Range("N136").Value = -23331.693
